# Start adding extra obs data
# Insert a new data row at row 21 of the CottonObserved sheet (this pushes
# the existing rows 21-39 down to 22-40) and populate it with the first
# readings for the new "ForestHill2023IrrigationFull" treatment at the
# 2023-12-20 (serial 45355) sampling date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CottonObserved")

# --- Insert the new row, shifting everything below it down by one ------
$ws.Rows(21).Insert()

# --- Populate the new row 21 -------------------------------------------
$ws.Cells.Item(21, 1).Value = "ForestHill2023IrrigationFull"   # A21 SimulationName
$ws.Cells.Item(21, 2).Value = 45355                             # B21 Clock.Today
$ws.Cells.Item(21, 4).Value = 23                                # D21 Cotton.Leaf.NodeNumber
$ws.Cells.Item(21, 5).Value = 79                                # E21 EMCalculator.Script.EMp100
$ws.Cells.Item(21, 6).Value = 5.51                              # F21 Cotton.Leaf.LAI
$ws.Cells.Item(21, 7).Value = 0.53                              # G21 Cotton.Leaf.LAIError

# --- Keep the filter database in sync with the extra row ---------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "CottonObserved!_FilterDatabase") {
        $n.RefersTo = "=CottonObserved!`$A`$1:`$EQ`$2581"
    }
}

# --- Re-apply the previous manual sort range, shifted down by one row --
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("B2:B2268")) | Out-Null
$sort.SortFields.Add($ws.Range("C2:C2268")) | Out-Null
$sort.SetRange($ws.Range("A1638:EQ1817"))
$sort.Header = 0
$sort.Apply()

# --- Move the active selection to the newly-added row -------------------
$ws.Range("A21").Select() | Out-Null

Write-Output "Inserted new observation row at row 21"
